$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9416471719741821
$ws.Range("B1").Value = 1.462187647819519
$ws.Range("C1").Value = 4.962259769439697
$ws.Range("D1").Value = 2.549798250198364
$ws.Range("E1").Value = 0.3758462369441986
